$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35; existing rows 35-68 shift down to 36-69.
$ws.Rows("35:35").Insert()

# Populate the newly inserted row 35 with the new weekly record.
$ws.Range("A35").Value = 4
$ws.Range("B35").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C35").Value = "Los Lagos"
$ws.Range("D35").Value = 44904
$ws.Range("E35").Value = 10
$ws.Range("F35").Value = 300000000
$ws.Range("G35").Value = "Espárragos"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 500
$ws.Range("K35").Value = 1500
$ws.Range("L35").Value = 1500
$ws.Range("M35").Value = 1500
$ws.Range("N35").Value = "$/kilo"
$ws.Range("O35").Value = "Provincia de Linares"
$ws.Range("P35").Value = 1500
$ws.Range("Q35").Value = 1
$ws.Range("R35").Value = "Hortaliza"
